$d = $word.ActiveDocument

# 1. Replace mill name in title
$d.Content.Find.Execute("HIDECO SUGAR MILLING COMPANY, INC.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Cotabato Sugar Central Company, Inc.", 2)

# 2. Replace address text (note the leading " of " prefix is removed from this run's text)
$d.Content.Find.Execute(" of 19/F Citibank Tower, Paseo De Roxas, Makati City", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "6/F Filinvest Bldg., No. 79 EDSA, Highway Hills, Mandaluyong City", 2)

# 3. Replace the day number
$d.Content.Find.Execute("Given this 10", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Given this 18", 2)

# 4. Replace license number suffix
$d.Content.Find.Execute("2021-05", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2021-01", 2)
